$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.093.53'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '2.055.33'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '247.41'
$ws.Range("E5").Value = '  -0.05%  '
$ws.Range("D6").Value = '0.663'
$ws.Range("E6").Value = '  +0.36%  '
$ws.Range("D7").Value = '57.71'
$ws.Range("E7").Value = '  +3.36%  '
$ws.Range("D9").Value = '0.381'
$ws.Range("E9").Value = '  +0.95%  '
$ws.Range("D10").Value = '0.0779'
$ws.Range("E10").Value = '  +0.24%  '
$ws.Range("D12").Value = '15.75'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '0.908'
$ws.Range("E13").Value = '  +15.15%  '
$ws.Range("D14").Value = '2.356.89'
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").Value = '5.79'
$ws.Range("E15").Value = '  +3.82%  '
$ws.Range("D16").Value = '2.061.27'
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").Value = '18.57'
$ws.Range("E17").Value = '  +13.86%  '
$ws.Range("D18").Value = '37.066.17'
$ws.Range("E18").Value = '  +0.34%  '
$ws.Range("D19").Value = '74.88'
$ws.Range("E19").Value = '  +1.54%  '
$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  +1.47%  '
$ws.Range("D21").Value = '5.49'
$ws.Range("E21").Value = '  +3.72%  '
$ws.Range("D22").Value = '237.25'
$ws.Range("E22").Value = '  +0.79%  '
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("E24").Value = '  +4.92%  '
$ws.Range("D25").Value = '9.59'
$ws.Range("E25").Value = '  +6.34%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").Value = '170.52'
$ws.Range("E26").Value = '  +1.90%  '
$ws.Range("B27").Value = 'PancakeSwap'
$ws.Range("C27").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D27").Value = '2.17'
$ws.Range("E27").Value = '  -0.21%  '
$ws.Range("D28").Value = '20.11'
$ws.Range("E28").Value = '  +1.98%  '
$ws.Range("D29").Value = '5.48'
$ws.Range("E29").Value = '  +18.31%  '
$ws.Range("D30").Value = '0.124'
$ws.Range("E30").Value = '  +0.45%  '
$ws.Range("D31").Value = '1.14'
$ws.Range("E31").Value = '  +3.57%  '
$ws.Range("D32").Value = '4.84'
$ws.Range("E32").Value = '  +10.40%  '
$ws.Range("D33").Value = '0.0621'
$ws.Range("E33").Value = '  +1.73%  '
$ws.Range("D34").Value = '0.0878'
$ws.Range("E34").Value = '  +1.01%  '
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = '  -0.17%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").Value = '2.31'
$ws.Range("E36").Value = '  +4.36%  '
$ws.Range("E37").Value = '  +4.52%  '
$ws.Range("D38").Value = '1.32'
$ws.Range("E38").Value = '  -0.75%  '
$ws.Range("D39").Value = '5.17'
$ws.Range("E39").Value = '  +6.13%  '
$ws.Range("D40").Value = '3.11'
$ws.Range("E40").Value = '  -2.16%  '
$ws.Range("D41").Value = "'0.100"
$ws.Range("E41").Value = '  -5.07%  '
$ws.Range("D42").Value = '0.0224'
$ws.Range("E42").Value = '  +1.84%  '
$ws.Range("D43").Value = '1.16'
$ws.Range("E43").Value = '  +4.81%  '
$ws.Range("D44").Value = "'99.30"
$ws.Range("E44").Value = '  +4.48%  '
$ws.Range("D45").Value = "'17.20"
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = "'2.40"
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("D47").Value = '1.303.16'
$ws.Range("E47").Value = '  +2.26%  '
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '2.87'
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").Value = '3.75'
$ws.Range("E49").Value = '  +10.83%  '
$ws.Range("D50").Value = "'6.90"
$ws.Range("E50").Value = '  +3.72%  '
$ws.Range("D51").Value = '2.242.36'
$ws.Range("E51").Value = '  +0.71%  '
